$d = $word.ActiveDocument

$d.Content.Find.Execute("923÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "792÷8=", 2) | Out-Null
$d.Content.Find.Execute("424÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "223÷3=", 2) | Out-Null
$d.Content.Find.Execute("176÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "511÷9=", 2) | Out-Null
$d.Content.Find.Execute("548÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "567÷7=", 2) | Out-Null
$d.Content.Find.Execute("995÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "451÷2=", 2) | Out-Null
$d.Content.Find.Execute("646÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "534÷7=", 2) | Out-Null
$d.Content.Find.Execute("543÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "562÷8=", 2) | Out-Null
$d.Content.Find.Execute("773÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "347÷7=", 2) | Out-Null
$d.Content.Find.Execute("188÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "356÷9=", 2) | Out-Null
$d.Content.Find.Execute("921÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "354÷5=", 2) | Out-Null
$d.Content.Find.Execute("287÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "973÷7=", 2) | Out-Null
$d.Content.Find.Execute("374÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "317÷5=", 2) | Out-Null
$d.Content.Find.Execute("240÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "264÷2=", 2) | Out-Null
$d.Content.Find.Execute("948÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "756÷3=", 2) | Out-Null
$d.Content.Find.Execute("691÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "292÷7=", 2) | Out-Null
$d.Content.Find.Execute("150÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "663÷8=", 2) | Out-Null
$d.Content.Find.Execute("849÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "347÷4=", 2) | Out-Null
$d.Content.Find.Execute("559÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "584÷3=", 2) | Out-Null
$d.Content.Find.Execute("378÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "309÷4=", 2) | Out-Null
$d.Content.Find.Execute("643÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "663÷6=", 2) | Out-Null
$d.Content.Find.Execute("779÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "431÷4=", 2) | Out-Null
$d.Content.Find.Execute("110÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "260÷6=", 2) | Out-Null
$d.Content.Find.Execute("415÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "167÷3=", 2) | Out-Null
$d.Content.Find.Execute("491÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "833÷6=", 2) | Out-Null
$d.Content.Find.Execute("528÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "765÷6=", 2) | Out-Null
